# Ubound fixed for L1&L2 Data
# Clamp column C (l2) values that exceed the valid UBound of 99 down to 99.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(17, 23, 37, 45, 52, 63, 93, 98, 115, 128, 133, 137, 140, 142, 163, 168, 172, 192, 198, 203, 207, 212)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = 99
}
